# New docking path: add a "shift left 10 cm" offset row and the four
# corresponding shifted points (Startpunkt / Zwischenwert / Einfahrt /
# Endpunkt) underneath the existing table on Tabelle1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# xlCenter, used below to reproduce the "vertical center" alignment that
# every style in this sheet already carries.
$xlCenter = -4108

# --- New labels, in the same order Excel originally appended them to the
#     shared-string table (B16..B19 first, then B14) -------------------
$ws1.Cells.Item(16, 2).Value = "Startpunkt"
$ws1.Cells.Item(17, 2).Value = "Zwischenwert"
$ws1.Cells.Item(18, 2).Value = "Einfahrt"
$ws1.Cells.Item(19, 2).Value = "Endpunkt"
$ws1.Cells.Item(14, 2).Value = "Verschiebung nach links 10 cm"

# --- Row 14: 10 cm (0.1 m, already-normalised units) left-shift vector --
$ws1.Cells.Item(14, 3).Formula = "=-D12*E14"
$ws1.Cells.Item(14, 4).Formula = "=C12*E14"
$ws1.Cells.Item(14, 5).Value   = 0.1

# --- Rows 16-19: original points translated by the shift vector --------
$ws1.Cells.Item(16, 3).Formula = "=C5+C`$14"
$ws1.Cells.Item(16, 4).Formula = "=D5+D`$14"

$ws1.Cells.Item(17, 3).Formula = "=C6+C`$14"
$ws1.Cells.Item(17, 4).Formula = "=D6+D`$14"

$ws1.Cells.Item(18, 3).Formula = "=C7+C`$14"
$ws1.Cells.Item(18, 4).Formula = "=D7+D`$14"

$ws1.Cells.Item(19, 3).Formula = "=C8+C`$14"
$ws1.Cells.Item(19, 4).Formula = "=D8+D`$14"

# --- Formatting: reproduce the existing styles used elsewhere in the
#     sheet (vertical-centred text labels in column B, "0.00" numbers
#     with vertical-centred alignment in C/D, "0.000" in E) ------------
$labelCells = @(
    $ws1.Cells.Item(14, 2),
    $ws1.Cells.Item(16, 2),
    $ws1.Cells.Item(17, 2),
    $ws1.Cells.Item(18, 2),
    $ws1.Cells.Item(19, 2)
)
foreach ($c in $labelCells) {
    $c.VerticalAlignment = $xlCenter
}

$numberCells = @(
    $ws1.Cells.Item(14, 3), $ws1.Cells.Item(14, 4),
    $ws1.Cells.Item(16, 3), $ws1.Cells.Item(16, 4),
    $ws1.Cells.Item(17, 3), $ws1.Cells.Item(17, 4),
    $ws1.Cells.Item(18, 3), $ws1.Cells.Item(18, 4),
    $ws1.Cells.Item(19, 3), $ws1.Cells.Item(19, 4)
)
foreach ($c in $numberCells) {
    $c.NumberFormat = "0.00"
    $c.VerticalAlignment = $xlCenter
}

$ws1.Cells.Item(14, 5).NumberFormat = "0.000"
$ws1.Cells.Item(14, 5).VerticalAlignment = $xlCenter

# --- New selection left behind after editing cell F26 -------------------
$ws1.Range("F26").Select()
